# Timesheet update for "DAY 6 (09-04-2022)" sheet.
# Aravindhan Ra was Absent; everyone else's reported activity for the day
# moved up/along one slot as reflected in the source diff.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("DAY 6 (09-04-2022)")
$ws.Activate()

# --- Row 8 (Aravindhan Ra) -> Absent ---
$ws.Range("C8").Value = "Absent"
$ws.Range("D8").Value = "Absent"
$ws.Range("F8").Value = 0
$ws.Range("G8").Value = 0

# --- Row 9 (Darshana) ---
$ws.Range("C9").Value = "Data  Model"
$ws.Range("D9").Value = "Brainstorming -30 Mins                                                  College Project review- 1 hour                                      Prototype and Acceptance criteria for TAC- 1hr                "
$ws.Range("E9").Value = "-"
$ws.Range("F9").Value = 1
$ws.Range("G9").Value = 1.5

# --- Row 10 (Deepika) ---
$ws.Range("C10").Value = "Data Model"
$ws.Range("D10").Value = "30 mins : Brainstorming with team`n1 Hr        : Prototype for Dashboard in TAC`n1 Hr        : Refining the TAC prototype (Alignment,Naming)                                                                                                                                    "
$ws.Range("E10").Value = "-"
$ws.Range("F10").Value = 2.5

# --- Row 11 (Gokul) ---
$ws.Range("C11").Value = "Data Model"
$ws.Range("D11").Value = "30 min: Brain Storming with team`n1.5 hr: Refined management Prototype(Dashboard,Profile)`n30 min: Refined management prototype(alignment)"
$ws.Range("E11").Value = "-"

# --- Row 12 (Kumaresh) ---
$ws.Range("C12").Value = "Data Model"
$ws.Range("D12").Value = "30 min: Brain Storming with team                                  1.5 hr: Refined management Prototype`n30 min: Refined management and interviewer prototype(alignment)"

# --- Row 13 (Prithvi) ---
$ws.Range("C13").Value = "Dependencies for Interviewer (7-13),Data Modelling "
$ws.Range("D13").Value = "30 mins : Brain stroming with Team`n1 hr : Refining the TAC Prototype ( Dash board , Profile )`n30 mins : TAC Prototype Page Alignments               30 mins : Refined TAC Prototype (Home page)  "
$ws.Range("F13").Value = 2.5
$ws.Range("G13").Value = 0

# --- Row 14 (Remuki) ---
$ws.Range("D14").Value = " 30 min :BrainStroming                                                      1 hr : Dependencies for TAC                                          1 hr :  Prototype for TAC"
$ws.Range("E14").Value = "-"
$ws.Range("G14").Value = 0
# E14 picks up the same "Status" cell formatting (centered, full border) as its neighbours
$ws.Range("E13").Copy()
$ws.Range("E14").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# --- Row 15 (Sheik Fareeth) ---
$ws.Range("D15").Value = "30 mins : Brain stroming`n1 hr : Refining the Interviewers Prototype ( Dash board , Profile , Navigations )`n1 hr : TAC ( Creating pools , Pool details prototype , and refining Slide )"
$ws.Range("G15").Value = 0

# --- Row 17 (Vishnu Prakaash R) ---
$ws.Range("C17").Value = "Exploring on Data models "
$ws.Range("D17").Value = "30 mins - Brainstorming`n1.5 hour - College Project review preparation and PPT`n1 hour - Refined Prototype for ADMIN   "
$ws.Range("E17").Value = "-"
# E17 loses its special bottom-border-less style and matches the rest of the column
$ws.Range("E16").Copy()
$ws.Range("E17").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# --- Row heights (auto-fit heights recorded by Excel after the text edits) ---
$ws.Rows.Item(8).RowHeight = 21
$ws.Rows.Item(9).RowHeight = 84
$ws.Rows.Item(10).RowHeight = 105
$ws.Rows.Item(14).RowHeight = 63
$ws.Rows.Item(16).RowHeight = 189

# --- View / selection: scrolled back up to the top, D3 selected ---
$ws.Range("D3").Select()
